# Auto-generated edit script applying the cryptos.xlsx diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "60.748.75"
$ws.Range("E2").Value = "  -4.42%  "
$ws.Range("D3").Value = "2.453.23"
$ws.Range("E3").Value = "  -5.88%  "
$ws.Range("E4").Value = "  -0.05%  "
$ws.Range("D5").Value = "'547.49"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -4.20%  "
$ws.Range("D6").Value = "'144.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -6.72%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -4.32%  "
$ws.Range("D9").Value = "2.447.84"
$ws.Range("E9").Value = "  -6.00%  "
$ws.Range("D10").Value = "'0.107"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -8.44%  "
$ws.Range("E11").Value = "  -2.00%  "
$ws.Range("D12").Value = "'5.36"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -7.87%  "
$ws.Range("E13").Value = "  -7.65%  "
$ws.Range("D14").Value = "'25.95"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -7.11%  "
$ws.Range("D15").Value = "2.891.97"
$ws.Range("E15").Value = "  -5.96%  "
$ws.Range("D16").Value = "'0.0000163"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -9.28%  "
$ws.Range("D17").Value = "60.687.51"
$ws.Range("E17").Value = "  -4.41%  "
$ws.Range("D18").Value = "2.455.98"
$ws.Range("E18").Value = "  -5.39%  "
$ws.Range("D19").Value = "'11.02"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -7.57%  "
$ws.Range("D20").Value = "'6.90"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -7.78%  "
$ws.Range("D21").Value = "'4.16"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -7.37%  "
$ws.Range("D22").Value = "'317.69"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -6.70%  "
$ws.Range("E23").Value = "  +0.04%  "
$ws.Range("D24").Value = "'63.29"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.63%  "
$ws.Range("D25").Value = "'1.77"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.73%  "
$ws.Range("B26").Value = "PEPE"
$ws.Range("C26").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D26").Value = "0.0₃0973"
$ws.Range("E26").Value = "  -8.62%  "
$ws.Range("B27").Value = "WrappedeETH"
$ws.Range("C27").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D27").Value = "2.577.20"
$ws.Range("E27").Value = "  -5.70%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -0.02%  "
$ws.Range("D29").Value = "'534.58"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -7.42%  "
$ws.Range("D30").Value = "'1.48"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -4.62%  "
$ws.Range("D31").Value = "'8.31"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -8.48%  "
$ws.Range("D32").Value = "'7.55"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -3.22%  "
$ws.Range("E33").Value = "  -7.27%  "
$ws.Range("E34").Value = "  -7.91%  "
$ws.Range("E35").Value = "  -8.84%  "
$ws.Range("D36").Value = "'5.82"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -11.51%  "
$ws.Range("D37").Value = "'0.998"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.15%  "
$ws.Range("D38").Value = "'4.81"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -9.72%  "
$ws.Range("E39").Value = "  -6.25%  "
$ws.Range("D40").Value = "'18.43"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -6.02%  "
$ws.Range("D41").Value = "'145.27"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -5.51%  "
$ws.Range("D43").Value = "'1.70"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -8.95%  "
$ws.Range("D44").Value = "'39.81"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.18%  "
$ws.Range("D45").Value = "'2.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.11%  "
$ws.Range("D46").Value = "'146.09"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -6.88%  "
$ws.Range("E47").Value = "  -7.92%  "
$ws.Range("D48").Value = "'20.76"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -11.54%  "
$ws.Range("D49").Value = "'0.0527"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -9.14%  "
$ws.Range("E50").Value = "  -7.53%  "
$ws.Range("E51").Value = "  -6.30%  "
